$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5250.278
$ws.Range("J40").Value = 7830
$ws.Range("L40").Value = 7830
$ws.Range("N40").Value = -8180
$ws.Range("H53").Value = 1220.1666
$ws.Range("I53").Value = 706
$ws.Range("J53").Value = 1734.3334
$ws.Range("K53").Value = 706
$ws.Range("L53").Value = 1734.3334
$ws.Range("M53").Value = -69
$ws.Range("N53").Value = -3008.3334
$ws.Range("H55").Value = 42
$ws.Range("I55").Value = 59.6
$ws.Range("K55").Value = 59.6
$ws.Range("M55").Value = 154.4
$ws.Range("H62").Value = 8297.6
$ws.Range("I62").Value = 7500
$ws.Range("J62").Value = 8497
$ws.Range("K62").Value = 7500
$ws.Range("L62").Value = 8497
$ws.Range("M62").Value = -6876
$ws.Range("N62").Value = -9745
$ws.Range("H65").Value = 8297.6
$ws.Range("I65").Value = 7500
$ws.Range("J65").Value = 8497
$ws.Range("K65").Value = 37500
$ws.Range("L65").Value = 42485
$ws.Range("M65").Value = -34380
$ws.Range("N65").Value = -48725
$ws.Range("H70").Value = 2604.1177
$ws.Range("I70").Value = 1012
$ws.Range("J70").Value = 3094
$ws.Range("K70").Value = 3036
$ws.Range("L70").Value = 9282
$ws.Range("M70").Value = -2766
$ws.Range("N70").Value = -9822
$ws.Range("H73").Value = 2604.1177
$ws.Range("I73").Value = 1012
$ws.Range("J73").Value = 3094
$ws.Range("K73").Value = 3036
$ws.Range("L73").Value = 9282
$ws.Range("M73").Value = -2100
$ws.Range("N73").Value = -11154
$ws.Range("H76").Value = 4898
$ws.Range("I76").Value = 4898
$ws.Range("K76").Value = 4898
$ws.Range("M76").Value = -4583
$ws.Range("H79").Value = 4898
$ws.Range("I79").Value = 4898
$ws.Range("K79").Value = 4898
$ws.Range("M79").Value = -3806
$ws.Range("H87").Value = 91284
$ws.Range("J87").Value = 91284
$ws.Range("L87").Value = 91284
$ws.Range("N87").Value = -93780
$ws.Range("H90").Value = 91284
$ws.Range("J90").Value = 91284
$ws.Range("L90").Value = 273852
$ws.Range("N90").Value = -286332
$ws.Range("H131").Value = 956.5
$ws.Range("J131").Value = 195
$ws.Range("L131").Value = 585
$ws.Range("N131").Value = -10665
$ws.Range("H132").Value = 2364.6538
$ws.Range("I132").Value = 1459.28
$ws.Range("K132").Value = 4377.84
$ws.Range("M132").Value = -1847.84
$ws.Range("H138").Value = 2015.1538
$ws.Range("I138").Value = 1318.4
$ws.Range("J138").Value = 4337.6665
$ws.Range("K138").Value = 3955.2
$ws.Range("L138").Value = 13012.9995
$ws.Range("M138").Value = 1184.8
$ws.Range("N138").Value = -23292.9995

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 5064.2856
$ws.Range("I45").Value = 2950
$ws.Range("J45").Value = 5416.6665
$ws.Range("K45").Value = 2950
$ws.Range("L45").Value = 5416.6665
$ws.Range("M45").Value = -2573
$ws.Range("N45").Value = -6170.6665
$ws.Range("H74").Value = 2286.8928
$ws.Range("I74").Value = 2040.1923
$ws.Range("J74").Value = 5494
$ws.Range("K74").Value = 2040.1923
$ws.Range("L74").Value = 5494
$ws.Range("M74").Value = -1166.1923
$ws.Range("N74").Value = -7242
$ws.Range("H77").Value = 2286.8928
$ws.Range("I77").Value = 2040.1923
$ws.Range("J77").Value = 5494
$ws.Range("K77").Value = 10200.9615
$ws.Range("L77").Value = 27470
$ws.Range("M77").Value = -5832.961499999999
$ws.Range("N77").Value = -36206
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 8489.9
$ws.Range("I86").Value = 7550
$ws.Range("J86").Value = 9116.5
$ws.Range("K86").Value = 7550
$ws.Range("L86").Value = 9116.5
$ws.Range("M86").Value = -6427
$ws.Range("N86").Value = -11362.5
$ws.Range("H89").Value = 8489.9
$ws.Range("I89").Value = 7550
$ws.Range("J89").Value = 9116.5
$ws.Range("K89").Value = 37750
$ws.Range("L89").Value = 45582.5
$ws.Range("M89").Value = -32134
$ws.Range("N89").Value = -56814.5
$ws.Range("H105").Value = 2807.375
$ws.Range("I105").Value = 2502.25
$ws.Range("K105").Value = 2502.25
$ws.Range("M105").Value = -755.25

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6772.9536
$ws.Range("I31").Value = 4330.8647
$ws.Range("K31").Value = 4330.8647
$ws.Range("M31").Value = -4035.8647
$ws.Range("H34").Value = 6772.9536
$ws.Range("I34").Value = 4330.8647
$ws.Range("K34").Value = 4330.8647
$ws.Range("M34").Value = -4128.8647
$ws.Range("H107").Value = 411.1905
$ws.Range("I107").Value = 355.27777
$ws.Range("J107").Value = 746.6667
$ws.Range("K107").Value = 355.27777
$ws.Range("L107").Value = 746.6667
$ws.Range("M107").Value = 1564.72223
$ws.Range("N107").Value = -4586.6667
$ws.Range("H134").Value = 2043.7632
$ws.Range("I134").Value = 1827.0286
$ws.Range("K134").Value = 5481.085800000001
$ws.Range("M134").Value = -2946.085800000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 40.588234
$ws.Range("J2").Value = 39.444443
$ws.Range("L2").Value = 236.666658
$ws.Range("N2").Value = -462.666658
$ws.Range("H38").Value = 299.27274
$ws.Range("J38").Value = 30
$ws.Range("L38").Value = 90
$ws.Range("N38").Value = -784
$ws.Range("H86").Value = 373.2857
$ws.Range("J86").Value = 475.375
$ws.Range("L86").Value = 1426.125
$ws.Range("N86").Value = -3798.125
$ws.Range("H89").Value = 373.2857
$ws.Range("J89").Value = 475.375
$ws.Range("L89").Value = 4278.375
$ws.Range("N89").Value = -16134.375

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2838.6
$ws.Range("I80").Value = 2673.25
$ws.Range("J80").Value = 3500
$ws.Range("K80").Value = 2673.25
$ws.Range("L80").Value = 3500
$ws.Range("M80").Value = -1675.25
$ws.Range("N80").Value = -5496
$ws.Range("H83").Value = 2838.6
$ws.Range("I83").Value = 2673.25
$ws.Range("J83").Value = 3500
$ws.Range("K83").Value = 13366.25
$ws.Range("L83").Value = 17500
$ws.Range("M83").Value = -8374.25
$ws.Range("N83").Value = -27484
$ws.Range("H97").Value = 1533.2858
$ws.Range("I97").Value = 1447.6
$ws.Range("J97").Value = 1747.5
$ws.Range("K97").Value = 1447.6
$ws.Range("L97").Value = 1747.5
$ws.Range("M97").Value = -951.5999999999999
$ws.Range("N97").Value = -2739.5
$ws.Range("H122").Value = 2925
$ws.Range("I122").Value = 1770.5
$ws.Range("J122").Value = 6388.5
$ws.Range("K122").Value = 5311.5
$ws.Range("L122").Value = 19165.5
$ws.Range("M122").Value = -2861.5
$ws.Range("N122").Value = -24065.5
$ws.Range("H126").Value = 3419
$ws.Range("I126").Value = 3419
$ws.Range("K126").Value = 10257
$ws.Range("M126").Value = -7787
$ws.Range("H128").Value = 150000
$ws.Range("J128").Value = 150000
$ws.Range("L128").Value = 150000
$ws.Range("N128").Value = -159960
$ws.Range("H129").Value = 45833.168
$ws.Range("J129").Value = 45833.168
$ws.Range("L129").Value = 45833.168
$ws.Range("N129").Value = -55833.168
$ws.Range("H130").Value = 49999
$ws.Range("J130").Value = 49999
$ws.Range("L130").Value = 49999
$ws.Range("N130").Value = -60039
$ws.Range("H131").Value = 99000
$ws.Range("J131").Value = 99000
$ws.Range("L131").Value = 99000
$ws.Range("N131").Value = -109080
$ws.Range("H132").Value = 66680.836
$ws.Range("I132").Value = 73453.56
$ws.Range("K132").Value = 220360.68
$ws.Range("M132").Value = -217830.68

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7842.8696
$ws.Range("J40").Value = 10277
$ws.Range("L40").Value = 10277
$ws.Range("N40").Value = -10549
$ws.Range("H61").Value = 5140.1875
$ws.Range("I61").Value = 4082.5557
$ws.Range("K61").Value = 4082.5557
$ws.Range("M61").Value = -3880.5557
$ws.Range("H105").Value = 38028.25
$ws.Range("J105").Value = 38028.25
$ws.Range("L105").Value = 38028.25
$ws.Range("N105").Value = -45016.25
$ws.Range("H113").Value = 5140.1875
$ws.Range("I113").Value = 4082.5557
$ws.Range("K113").Value = 4082.5557
$ws.Range("M113").Value = -1912.5557
$ws.Range("H132").Value = 9874.308000000001
$ws.Range("I132").Value = 8249.166999999999
$ws.Range("K132").Value = 24747.501
$ws.Range("M132").Value = -22217.501

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 11333.223
$ws.Range("I62").Value = 10333
$ws.Range("K62").Value = 10333
$ws.Range("M62").Value = -9709
$ws.Range("H65").Value = 11333.223
$ws.Range("I65").Value = 10333
$ws.Range("K65").Value = 51665
$ws.Range("M65").Value = -48545
$ws.Range("H113").Value = 681.3333
$ws.Range("I113").Value = 613.55554
$ws.Range("J113").Value = 783
$ws.Range("K113").Value = 1840.66662
$ws.Range("L113").Value = 2349
$ws.Range("M113").Value = 329.33338
$ws.Range("N113").Value = -6689
$ws.Range("H126").Value = 5663.636
$ws.Range("I126").Value = 2890
$ws.Range("K126").Value = 8670
$ws.Range("M126").Value = -6200
